$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting rows 192:281 down to 193:282.
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with its data.
$ws.Cells.Item(192, 1).Value = 3
$ws.Cells.Item(192, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 45141
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 100112026
$ws.Cells.Item(192, 7).Value = "Haba"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 100
$ws.Cells.Item(192, 11).Value = 14500
$ws.Cells.Item(192, 12).Value = 15000
$ws.Cells.Item(192, 13).Value = 14750
$ws.Cells.Item(192, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(192, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(192, 16).Value = 590
$ws.Cells.Item(192, 17).Value = 25
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Apply the same date number format used by the other rows in column D.
$ws.Cells.Item(192, 4).NumberFormat = $ws.Cells.Item(193, 4).NumberFormat
